# Apply the data edits from the commit "Add files via upload":
#   - I9 changes from "S4_NUM" to "BAD" (now matching the rest of row 9)
#   - Row 26 (B26:AA26) changes from "SG" to "BAD" (now matching rows 25/27)
# and mirror the author's last-saved view state (zoom level + active cell)
# as closely as the object model allows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Cell value corrections ---------------------------------------------
$ws.Range("I9").Value = "BAD"
$ws.Range("B26:AA26").Value = "BAD"

# --- View state: zoom + scroll position + active cell -------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 4   # best-effort: corresponds to topLeftCell="D1"
$win.ScrollRow = 1
$win.Zoom = 150
$ws.Range("I9").Select() | Out-Null
